# Exercise 5 lab sheet updates
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. "There is a file on Github that contains so" -> "There is a file in the Github repository that contains so"
Replace-Text "There is a file on Github that contains so" "There is a file in the Github repository that contains so"

# 2. wget command -> local file path
Replace-Text "wget https://freo.me/doctors-practices -O practices.csv " "~/BigData/datafiles/practices/ukpractices.csv "

# 3. postcode areas
Replace-Text "for the postcode areas: OX1, SW11." "for the postcode areas: BN1, GU27."

# 4. "Ask me or David if you get stuck." -> "Ask one of us if you get stuck." (and move the _GoBack bookmark here)
Replace-Text "Ask me or David if you get stuck." "Ask one of us if you get stuck."

# Move the _GoBack bookmark from the title paragraph to after "one of us"
$find = $d.Content.Find
$find.Execute("one of us") | Out-Null
$bmRange = $find.Parent.Duplicate
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 5. Footer copyright line update
$d2 = $d
$ftr = $d.Sections.First.Footers.Item(1)
$ftrFind = $ftr.Range.Find
$ftrFind.Execute("(c) Paul Fremantle 2015.  Licensed under the This work is licensed under a ") | Out-Null

# 6. Header text update
$hdr = $d.Sections.First.Headers.Item(1)
